# "case with 380 kV done" - update loading_percent values for rows 2-25
# (columns B, D, E, F, G, H, I, K, O); columns A, C, J, L, M, N are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each triple is (row, column, new value)
$updates = @(
    @(2, 2, 6.677490187444684),
    @(2, 4, 3.837303855917068),
    @(2, 5, 16.58317341930333),
    @(2, 6, 18.72891597484025),
    @(2, 7, 19.32121463004367),
    @(2, 8, 11.81178719602177),
    @(2, 9, 17.95775239994621),
    @(2, 11, 12.41835138854054),
    @(2, 15, 16.71993153412754),
    @(3, 2, 6.516664457825418),
    @(3, 4, 3.753607368066372),
    @(3, 5, 15.63784156558004),
    @(3, 6, 18.74085985622836),
    @(3, 7, 19.31289239607152),
    @(3, 8, 11.86246803330413),
    @(3, 9, 18.07193712749438),
    @(3, 11, 11.89613629296124),
    @(3, 15, 16.79104181059758),
    @(4, 2, 6.416440836923975),
    @(4, 4, 3.700721285431267),
    @(4, 5, 15.03195387295727),
    @(4, 6, 18.75532428746634),
    @(4, 7, 19.31845628124648),
    @(4, 8, 11.8960772284795),
    @(4, 9, 18.14599054324907),
    @(4, 11, 11.56172679841083),
    @(4, 15, 16.83993333260222),
    @(5, 2, 6.375283164120018),
    @(5, 4, 3.678811801350371),
    @(5, 5, 14.77891697539674),
    @(5, 6, 18.76300565788661),
    @(5, 7, 19.32339455225744),
    @(5, 8, 11.91039854198856),
    @(5, 9, 18.17716133424684),
    @(5, 11, 11.42212127673689),
    @(5, 15, 16.86116515558683),
    @(6, 2, 6.368431566706157),
    @(6, 4, 3.675152680827797),
    @(6, 5, 14.73653883627358),
    @(6, 6, 18.76438886584894),
    @(6, 7, 19.32437540783221),
    @(6, 8, 11.91281432008936),
    @(6, 9, 18.18239727405659),
    @(6, 11, 11.39874290258496),
    @(6, 15, 16.86476948722092),
    @(7, 2, 6.415886973333258),
    @(7, 4, 3.700427231049813),
    @(7, 5, 15.02856576850157),
    @(7, 6, 18.75542065486362),
    @(7, 7, 19.31851208620941),
    @(7, 8, 11.89626784037142),
    @(7, 9, 18.14640689870381),
    @(7, 11, 11.55985732955716),
    @(7, 15, 16.84021438507557),
    @(8, 2, 6.62237619819524),
    @(8, 4, 3.808766427423032),
    @(8, 5, 16.26264962500409),
    @(8, 6, 18.73155100026782),
    @(8, 7, 19.31612314805637),
    @(8, 8, 11.82874430181026),
    @(8, 9, 17.99630588536043),
    @(8, 11, 12.24122194264354),
    @(8, 15, 16.74336088660716),
    @(9, 2, 7.013301072452585),
    @(9, 4, 4.008603110282069),
    @(9, 5, 18.58489426020654),
    @(9, 6, 18.74152470894693),
    @(9, 7, 19.39651385105241),
    @(9, 8, 11.71614383585022),
    @(9, 9, 17.73316748759724),
    @(9, 11, 13.46328110769473),
    @(9, 15, 16.59522515730738),
    @(10, 2, 7.289133639242711),
    @(10, 4, 4.146813051990575),
    @(10, 5, 20.23478008578898),
    @(10, 6, 18.78363885939404),
    @(10, 7, 19.50767108228434),
    @(10, 8, 11.64556131491208),
    @(10, 9, 17.55875364878695),
    @(10, 11, 14.28617571097491),
    @(10, 15, 16.51226628770436),
    @(11, 2, 7.411609700989348),
    @(11, 4, 4.207646512138255),
    @(11, 5, 20.94282329481781),
    @(11, 6, 18.81035273388329),
    @(11, 7, 19.56950896105233),
    @(11, 8, 11.61610323765046),
    @(11, 9, 17.48349237826131),
    @(11, 11, 14.64338915665895),
    @(11, 15, 16.48022569232797),
    @(12, 2, 7.457516396498945),
    @(12, 4, 4.230376406652899),
    @(12, 5, 21.20486594294235),
    @(12, 6, 18.82155241334333),
    @(12, 7, 19.59453597505802),
    @(12, 8, 11.60533060112973),
    @(12, 9, 17.45557808455405),
    @(12, 11, 14.77613643665805),
    @(12, 15, 16.46891834222604),
    @(13, 2, 7.447651187951034),
    @(13, 4, 4.225494944272632),
    @(13, 5, 21.14870016183494),
    @(13, 6, 18.81909222694927),
    @(13, 7, 19.58907455240141),
    @(13, 8, 11.60763364716225),
    @(13, 9, 17.46156390536386),
    @(13, 11, 14.74765992616947),
    @(13, 15, 16.47131675920966),
    @(14, 2, 7.415396159964348),
    @(14, 4, 4.209522727207364),
    @(14, 5, 20.96450332282399),
    @(14, 6, 18.81125244135595),
    @(14, 7, 19.5715357503395),
    @(14, 8, 11.6152092925191),
    @(14, 9, 17.4811841248262),
    @(14, 11, 14.65436119786383),
    @(14, 15, 16.47927884543056),
    @(15, 2, 7.395576361470837),
    @(15, 4, 4.199699004501046),
    @(15, 5, 20.85088689426097),
    @(15, 6, 18.8065913644638),
    @(15, 7, 19.56100205682754),
    @(15, 8, 11.61989944807773),
    @(15, 9, 17.49327829128863),
    @(15, 11, 14.59688304588987),
    @(15, 15, 16.48426357206602),
    @(16, 2, 7.281065584860932),
    @(16, 4, 4.142795259076922),
    @(16, 5, 20.1876554719381),
    @(16, 6, 18.78204488855399),
    @(16, 7, 19.50385584906105),
    @(16, 8, 11.64753991851864),
    @(16, 9, 17.5637541508313),
    @(16, 11, 14.2624810881243),
    @(16, 15, 16.51447542313968),
    @(17, 2, 7.210018128482217),
    @(17, 4, 4.107354685804505),
    @(17, 5, 19.76992340524036),
    @(17, 6, 18.76892022672663),
    @(17, 7, 19.4716793854253),
    @(17, 8, 11.66517623654662),
    @(17, 9, 17.60803300220001),
    @(17, 11, 14.05290486020943),
    @(17, 15, 16.5344733964659),
    @(18, 2, 7.168873209009859),
    @(18, 4, 4.086779012193123),
    @(18, 5, 19.52565287248186),
    @(18, 6, 18.76208274873292),
    @(18, 7, 19.45423420725439),
    @(18, 8, 11.67556955171718),
    @(18, 9, 17.63388516911878),
    @(18, 11, 13.93075308455153),
    @(18, 15, 16.54651164159019),
    @(19, 2, 7.154895399666539),
    @(19, 4, 4.079780025486529),
    @(19, 5, 19.44225868683073),
    @(19, 6, 18.75988995090074),
    @(19, 7, 19.44851019746398),
    @(19, 8, 11.67913133719182),
    @(19, 9, 17.64270428466245),
    @(19, 11, 13.88912019589429),
    @(19, 15, 16.55067943347811),
    @(20, 2, 7.217610567218362),
    @(20, 4, 4.111147282265949),
    @(20, 5, 19.81480565740025),
    @(20, 6, 18.77024374866521),
    @(20, 7, 19.47499478760468),
    @(20, 8, 11.66327300128529),
    @(20, 9, 17.60327969479368),
    @(20, 11, 14.07538159953602),
    @(20, 15, 16.53228906308157),
    @(21, 2, 7.424883370336333),
    @(21, 4, 4.214222570643664),
    @(21, 5, 21.01877108954127),
    @(21, 6, 18.81352579584545),
    @(21, 7, 19.57664372244129),
    @(21, 8, 11.61297374980512),
    @(21, 9, 17.4754053096613),
    @(21, 11, 14.6818341569681),
    @(21, 15, 16.47691772817643),
    @(22, 2, 7.557578483415036),
    @(22, 4, 4.279796814715268),
    @(22, 5, 21.77022041731249),
    @(22, 6, 18.8481272959736),
    @(22, 7, 19.65245564189202),
    @(22, 8, 11.58233042523464),
    @(22, 9, 17.3952443095574),
    @(22, 11, 15.06346674700515),
    @(22, 15, 16.44554509261397),
    @(23, 2, 7.487022615496105),
    @(23, 4, 4.244966613989486),
    @(23, 5, 21.37238702119899),
    @(23, 6, 18.8290834589671),
    @(23, 7, 19.61113984881342),
    @(23, 8, 11.59848082333388),
    @(23, 9, 17.43771591369357),
    @(23, 11, 14.86114612991531),
    @(23, 15, 16.46184658273177),
    @(24, 2, 7.214178953215169),
    @(24, 4, 4.10943327263101),
    @(24, 5, 19.79452719853556),
    @(24, 6, 18.7696431785769),
    @(24, 7, 19.47349261188074),
    @(24, 8, 11.66413266293914),
    @(24, 9, 17.60542743249285),
    @(24, 11, 14.06522504295823),
    @(24, 15, 16.53327491499684),
    @(25, 2, 6.909350745355206),
    @(25, 4, 3.95599183782751),
    @(25, 5, 17.93940150751057),
    @(25, 6, 18.73272347265279),
    @(25, 7, 19.36561741493569),
    @(25, 8, 11.74447730780411),
    @(25, 9, 17.80102404223196),
    @(25, 11, 13.14550905988707),
    @(25, 15, 16.63078356446546),
)

foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}
